$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.965.72'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '2.365.26'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''303.47'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = '''95.34'
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '''0.501'
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("D10").Value = '''34.18'
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '''18.49'
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("D14").Value = '''6.74'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '2.730.91'
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '2.347.87'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").Value = '''0.796'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '42.945.24'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").Value = '''11.94'
$ws.Range("E19").Value = '  -2.25%  '
$ws.Range("D20").Value = '''6.27'
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '''235.00'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("E24").Value = '  -2.48%  '
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '''24.52'
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("E28").Value = '  +15.22%  '
$ws.Range("E29").Value = '  +2.37%  '
$ws.Range("D30").Value = '''32.07'
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").Value = '''17.54'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").Value = '''0.0717'
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("E35").Value = '  +3.46%  '
$ws.Range("E36").Value = '  +1.82%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '''4.34'
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '''2.84'
$ws.Range("E38").Value = '  +3.28%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '''123.63'
$ws.Range("E39").Value = '  -10.04%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '''2.26'
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").Value = '''21.18'
$ws.Range("E42").Value = '  -5.27%  '
$ws.Range("D43").Value = '1.933.61'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E45").Value = '  +4.19%  '
$ws.Range("D46").Value = '''9.30'
$ws.Range("E46").Value = '  -7.28%  '
$ws.Range("D47").Value = '''2.71'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").Value = '2.590.87'
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("D50").Value = '''71.48'
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("D51").Value = '''1.14'
$ws.Range("E51").Value = '  +0.79%  '
